# Atualização automática de SANTA_ROSA.xlsx
#
# 1. Rename "Paineis DARQ" -> "PAINEIS DARQ"
# 2. Rename "Recolhimento x Eliminacao" -> "RECOLHIMENTO X ELIMINAÇÃO"
# 3. Delete worksheet "Desarquivamentos Pendentes"

$wb = $excel.ActiveWorkbook

$wsPaineis = $wb.Worksheets.Item("Paineis DARQ")
$wsPaineis.Name = "PAINEIS DARQ"

$wsRecolhimento = $wb.Worksheets.Item("Recolhimento x Eliminacao")
$wsRecolhimento.Name = "RECOLHIMENTO X ELIMINAÇÃO"

$excel.DisplayAlerts = $false
$wsDesarquivamentos = $wb.Worksheets.Item("Desarquivamentos Pendentes")
$wsDesarquivamentos.Delete()
$excel.DisplayAlerts = $true
